$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in F1 from "Aperp" to "value"
$ws.Range("F1").Value = "value"

# Update the selection to match the target state (single cell F11)
$ws.Range("F11").Select()
